$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    "baseball pants high knee",
    "knee pack",
    "girls compression pants",
    "basketball gear for boys",
    "gym knee pads",
    "extra large knee pads",
    "tall leggings",
    "large knee pads",
    "football shock pads",
    "fitness compression pants men",
    "spandex leggings",
    "knee pad hiking",
    "workout legging for men",
    "compression legs men",
    "pants for men",
    "mens baseball clothing",
    "capri legging pants",
    "compression pants knee length",
    "knee pad for pain",
    "camo knee pads for basketball",
    "mcdavid basketball tights knee pads",
    "skateboard knee pads youth",
    "kids basketball knee pads",
    "youth basketball pants tearaway",
    "kneee pads for basketball",
    "commpression pants for basketball",
    "supreme basketball leggings",
    "volleyball knee pads youth girls",
    "adidas knee pads basketball",
    "basketball tights men",
    "compression pants men adidas",
    "pant with knee pads",
    "knee pads tights",
    "compression pants with knee",
    "women basketball pants",
    "basketball knee pads for youth kids",
    "pantalones con rodilleras",
    "caterpillar knee pad pants",
    "thick mens leggings",
    "capri leggings men",
    "compression leggings men basketball",
    "capri tights men",
    "mens compression knee pads",
    "legging men",
    "basketball knee pads youth",
    "hex pads basketball knee",
    "baseball knee high pants",
    "little boys baseball pants",
    "boys baseball pants",
    "knee pad compression",
    "black compression tights",
    "capris men",
    "hip pads for volleyball",
    "soccer pants men",
    "boys knee pads",
    "capris pants men",
    "knee pad baseball",
    "volleyball kneepads youth",
    "compression leggings for boys",
    "knee pads lightweight",
    "knee pads sliding",
    "capri leggings mesh",
    "men running pants",
    "knee pad black",
    "baseball pants youth boys",
    "patella protector",
    "leggings pack",
    "knee pads for workout",
    "capri pants boys",
    "gym pad men",
    "mens athletic pants",
    "leggings youth",
    "knee pads protection",
    "knee pad softball",
    "knee pads fitness",
    "leggings youth boys",
    "gel knee pad",
    "basketballs for boys",
    "women athletic leggings",
    "knee pads insert",
    "capri leggings",
    "elbow pads basketball youth",
    "nike compression tights",
    "nike compression leggings men",
    "gloves and knee pads",
    "black leggings running",
    "lupo compression leggings",
    "compression leggings 30",
    "compression leggings black",
    "compression leggings circulation",
    "compression leggings girls",
    "compression leggings men under armour",
    "compression leggings men white",
    "compression leggings men nike",
    "compression leggings running",
    "compression leggings xxl",
    "running capri pants",
    "knee basketball pads",
    "girls black leggings",
    "basketball knee pads kids youth"
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $newValues[$i]
}
